$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-added stimulus row (row 15) with its Q1/Q1_ans/Q2/Q2_ans data.
$ws.Range("D15").Value = "According to common wisdom, feeling close to your partner is not important."
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = "According to the new study, the degree to which a person wants to be close to their partner is what matters the most."
$ws.Range("G15").Value = 1

# Re-fill H2:H16 as one operation so Excel stores it as a shared formula
# (matches the t="shared" ref="H2:H16" si="0" pattern produced by filling down).
$ws.Range("H2:H16").Formula = '=IF(ISBLANK(C2),0,LEN(TRIM(C2))-LEN(SUBSTITUTE(C2," ",""))+1)'

# Move the active selection to F16.
$ws.Range("F16").Select()
